# Update column C (Förändrad) values from 45203 to 45204 for data rows 2-40
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 40; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45203) {
        $cell.Value2 = 45204
    }
}
